# Creating common Utility for runmodes
# Update the Suite sheet: CustomerSuite's Runmode flag changes from "N" to "Y",
# and the saved selection moves from B5 to B4.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Suite")

# Change B3 (CustomerSuite / Runmode) value from "N" to "Y"
$ws.Range("B3").Value = "Y"

# Update the active selection to B4
$ws.Range("B4").Select()
